$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update conversion text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.18 = 4420.06 pesos`n✅ 4420.06 pesos = 1.17 = 760.3 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 850
$wsTasas.Range("O10").Value = 3757.05
$wsTasas.Range("N12").Value = 3779.2
$wsTasas.Range("O12").Value = 650.061
